$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the test data row (TC001) with the new account credentials
$ws.Range("D9").Value = "sal1@gmail.com"
$ws.Range("E9").Value = "Salsa123!"

# Update the active selection to match the new state (cell E9)
$ws.Range("E9").Select()
